$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,24
$arr[0,0] = -5.8247149282878
$arr[0,1] = 14.3170224031493
$arr[0,2] = -0.406838430804336
$arr[0,3] = 0.691035229236975
$arr[0,4] = -6.77822766096045
$arr[0,5] = 17.1872188841083
$arr[0,6] = -0.39437605971422
$arr[0,7] = 0.700677542957408
$arr[0,8] = 12.560879405376
$arr[0,9] = 0.108229728091123
$arr[0,10] = 116.057571490898
$arr[0,11] = 0.000000000000157542293495218
$arr[0,12] = 12.9700432204871
$arr[0,13] = 0.126756380585497
$arr[0,14] = 102.322606251279
$arr[0,15] = 0.00000000000000000000000126321232334932
$arr[0,16] = 13.4075151575492
$arr[0,17] = 0.0486154986536355
$arr[0,18] = 275.786848409639
$arr[0,19] = 0.0000000000000000000000000000476306318542484
$arr[0,20] = 13.5188744266351
$arr[0,21] = 0.0602542226915907
$arr[0,22] = 224.363933725126
$arr[0,23] = 0.000000000000000000000000000306971284999252
$ws.Range("B2:Y2").Value = $arr

$arr = New-Object 'object[,]' 1,24
$arr[0,0] = 0.0184235660613105
$arr[0,1] = 0.011140273702873
$arr[0,2] = 1.65378037853407
$arr[0,3] = 0.118997870577675
$arr[0,4] = 0.0165643124624637
$arr[0,5] = 0.0115667829539916
$arr[0,6] = 1.43205872612553
$arr[0,7] = 0.172690068273722
$arr[0,8] = 0.0208611678885488
$arr[0,9] = 0.0121320784065653
$arr[0,10] = 1.71950486878322
$arr[0,11] = 0.106134941477079
$arr[0,12] = 0.0148906335520594
$arr[0,13] = 0.0125377084648194
$arr[0,14] = 1.18766787358649
$arr[0,15] = 0.253481927368807
$arr[0,16] = 0.0376992230617386
$arr[0,17] = 0.0245304195100736
$arr[0,18] = 1.53683564385261
$arr[0,19] = 0.145216783890227
$arr[0,20] = 0.0793657814536189
$arr[0,21] = 0.0239364393551958
$arr[0,22] = 3.31568869855287
$arr[0,23] = 0.00428512670172188
$ws.Range("B3:Y3").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.00010524037515345
$arr[0,1] = 0.0000725543322433429
$arr[0,2] = 1.45050435858855
$arr[0,3] = 0.175751200660615
$arr[0,4] = 0.000145773567874724
$arr[0,5] = 0.0000772278008418024
$arr[0,6] = 1.88757890663409
$arr[0,7] = 0.0864049339541426
$arr[0,8] = 0.00019495048681685
$arr[0,9] = 0.00011696301854198
$arr[0,10] = 1.66677031122345
$arr[0,11] = 0.124371406377972
$arr[0,12] = 0.000178599771047061
$arr[0,13] = 0.000156613896993262
$arr[0,14] = 1.14038265106669
$arr[0,15] = 0.279108843806418
$arr[0,16] = 0.00111210185542342
$arr[0,17] = 0.000262888758091524
$arr[0,18] = 4.23031347364136
$arr[0,19] = 0.00156817030100546
$ws.Range("B4:U4").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.0348435761308152
$arr[0,1] = 0.0111457845566277
$arr[0,2] = 3.12616630563668
$arr[0,3] = 0.00688582779501879
$arr[0,4] = 0.0388671226089921
$arr[0,5] = 0.0113105091375522
$arr[0,6] = 3.4363725042181
$arr[0,7] = 0.00363920748337796
$arr[0,8] = 0.0312848756291325
$arr[0,9] = 0.0139198119934989
$arr[0,10] = 2.24750705280674
$arr[0,11] = 0.0399601973859502
$arr[0,12] = 0.0335539207169659
$arr[0,13] = 0.0156696369555417
$arr[0,14] = 2.14133363856264
$arr[0,15] = 0.048958593866375
$arr[0,16] = 0.0673315457662609
$arr[0,17] = 0.0155630797931369
$arr[0,18] = 4.32636384708079
$arr[0,19] = 0.000590979573726712
$ws.Range("B5:U5").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.0425495439358445
$arr[0,1] = 0.0166840733449311
$arr[0,2] = 2.55030909156076
$arr[0,3] = 0.0221296846947093
$arr[0,4] = 0.0508034648033069
$arr[0,5] = 0.014660903965056
$arr[0,6] = 3.46523413047355
$arr[0,7] = 0.00346022381147689
$arr[0,8] = 0.0418095444581613
$arr[0,9] = 0.0166215787284148
$arr[0,10] = 2.51537745850142
$arr[0,11] = 0.0237600260946259
$arr[0,12] = 0.0456286219673965
$arr[0,13] = 0.0190510625530982
$arr[0,14] = 2.39506966292419
$arr[0,15] = 0.0301155720310971
$arr[0,16] = 0.0772280967180565
$arr[0,17] = 0.0160884233301552
$arr[0,18] = 4.80022778697679
$arr[0,19] = 0.00023471583302433
$ws.Range("B6:U6").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.0328778108518981
$arr[0,1] = 0.0177410401712181
$arr[0,2] = 1.85320649379042
$arr[0,3] = 0.0828637734398298
$arr[0,4] = 0.0382211942274417
$arr[0,5] = 0.0149972292175555
$arr[0,6] = 2.5485503804063
$arr[0,7] = 0.0218102994833571
$arr[0,8] = 0.0292068208892011
$arr[0,9] = 0.0160930849599316
$arr[0,10] = 1.81486774983913
$arr[0,11] = 0.0888805448670006
$arr[0,12] = 0.0393131106243669
$arr[0,13] = 0.0190882169243255
$arr[0,14] = 2.05954860950199
$arr[0,15] = 0.0565977508924102
$arr[0,16] = 0.0711513049424382
$arr[0,17] = 0.0144136203890769
$arr[0,18] = 4.9363937041355
$arr[0,19] = 0.000162295033362885
$ws.Range("B7:U7").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.0128032925231662
$arr[0,1] = 0.0234789357810852
$arr[0,2] = 0.545309746682839
$arr[0,3] = 0.59375229571626
$arr[0,4] = 0.0220143396979452
$arr[0,5] = 0.0197300253212228
$arr[0,6] = 1.11577858312555
$arr[0,7] = 0.282525885115434
$arr[0,8] = 0.0125083796406478
$arr[0,9] = 0.0196706685439126
$arr[0,10] = 0.635889909523115
$arr[0,11] = 0.534678535166805
$arr[0,12] = 0.0264518282225683
$arr[0,13] = 0.0264187262324018
$arr[0,14] = 1.00125297449526
$arr[0,15] = 0.332997150796834
$arr[0,16] = 0.0665973510369997
$arr[0,17] = 0.0330941260024397
$arr[0,18] = 2.01236168110589
$arr[0,19] = 0.0630005247558887
$ws.Range("B8:U8").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.080383264414336
$arr[0,1] = 0.0263803262604357
$arr[0,2] = 3.04709136728502
$arr[0,3] = 0.00807656712731205
$arr[0,4] = 0.091115845030034
$arr[0,5] = 0.0215933598113075
$arr[0,6] = 4.21962333913042
$arr[0,7] = 0.000729677342378604
$arr[0,8] = 0.0807708125164583
$arr[0,9] = 0.0212767305668269
$arr[0,10] = 3.79620413309131
$arr[0,11] = 0.00173205638073489
$arr[0,12] = 0.0907186651610932
$arr[0,13] = 0.0300944642890274
$arr[0,14] = 3.01446353355324
$arr[0,15] = 0.00864801584526214
$arr[0,16] = 0.117666500286746
$arr[0,17] = 0.0523836218994938
$arr[0,18] = 2.24624598338214
$arr[0,19] = 0.0400720259140853
$ws.Range("B9:U9").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.119490239414695
$arr[0,1] = 0.029373481639749
$arr[0,2] = 4.0679630995121
$arr[0,3] = 0.000984108218542663
$arr[0,4] = 0.128970731110046
$arr[0,5] = 0.0250520835218218
$arr[0,6] = 5.14810398894389
$arr[0,7] = 0.000114476989694538
$arr[0,8] = 0.106813189913608
$arr[0,9] = 0.0252320552993904
$arr[0,10] = 4.23323382285822
$arr[0,11] = 0.000704731153899889
$arr[0,12] = 0.119403870036219
$arr[0,13] = 0.0347539219457436
$arr[0,14] = 3.43569483244588
$arr[0,15] = 0.00362366676374271
$arr[0,16] = 0.128232774092908
$arr[0,17] = 0.0547831552329123
$arr[0,18] = 2.34073363514246
$arr[0,19] = 0.0333321394426924
$ws.Range("B10:U10").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = 0.0808663054811024
$arr[0,1] = 0.0370339449756134
$arr[0,2] = 2.18357254498143
$arr[0,3] = 0.0437037628932784
$arr[0,4] = 0.0823554225229931
$arr[0,5] = 0.0342276171877639
$arr[0,6] = 2.40611030768553
$arr[0,7] = 0.0281283658157955
$arr[0,8] = 0.0483739477476084
$arr[0,9] = 0.04113811076508
$arr[0,10] = 1.17589132918254
$arr[0,11] = 0.25627671088553
$arr[0,12] = 0.116415092903591
$arr[0,13] = 0.0465518319210378
$arr[0,14] = 2.50076287225509
$arr[0,15] = 0.023256354419546
$arr[0,16] = 0.442252569345002
$arr[0,17] = 0.0607209669426589
$arr[0,18] = 7.28335847751961
$arr[0,19] = 0.00000151807112591651
$ws.Range("B11:U11").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.0249792342488939
$arr[0,1] = 0.0227554931286625
$arr[0,2] = 1.09772326653868
$arr[0,3] = 0.286294429608753
$arr[0,4] = 0.0141705951065409
$arr[0,5] = 0.0244573599551447
$arr[0,6] = 0.579400030605514
$arr[0,7] = 0.569282414391197
$arr[0,8] = -0.0123024717404969
$arr[0,9] = 0.0237652318748445
$arr[0,10] = -0.517666808608717
$arr[0,11] = 0.610818961801791
$arr[0,12] = 0.0327024964726892
$arr[0,13] = 0.0185145387312341
$arr[0,14] = 1.76631440552824
$arr[0,15] = 0.0941589511948885
$ws.Range("B12:Q12").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.0537415865778203
$arr[0,1] = 0.0146484004409232
$arr[0,2] = 3.66876825866137
$arr[0,3] = 0.00170974496050321
$arr[0,4] = 0.0601655301544533
$arr[0,5] = 0.0171847280247423
$arr[0,6] = 3.50110458936725
$arr[0,7] = 0.00248874679001274
$arr[0,8] = 0.0905330419484441
$arr[0,9] = 0.0246517126328479
$arr[0,10] = 3.6724848815498
$arr[0,11] = 0.00170774595674166
$arr[0,12] = 0.0971638202630247
$arr[0,13] = 0.0265888802874225
$arr[0,14] = 3.65430282180731
$arr[0,15] = 0.00175279298698715
$ws.Range("B13:Q13").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.0584810147844914
$arr[0,1] = 0.00953282454732904
$arr[0,2] = 6.13469958396297
$arr[0,3] = 0.0000143857496844711
$arr[0,4] = 0.0587088295231157
$arr[0,5] = 0.0097737869511503
$arr[0,6] = 6.00676378731645
$arr[0,7] = 0.0000184086792331917
$arr[0,8] = 0.0496888490527145
$arr[0,9] = 0.0134338976528695
$arr[0,10] = 3.69876638460922
$arr[0,11] = 0.00196256123746919
$arr[0,12] = 0.0671237856922173
$arr[0,13] = 0.017988987353885
$arr[0,14] = 3.73138211572099
$arr[0,15] = 0.00185619339627581
$ws.Range("B14:Q14").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = -0.100391732854145
$arr[0,1] = 0.0296096087575088
$arr[0,2] = -3.39051196779783
$arr[0,3] = 0.00400366332627309
$arr[0,4] = -0.105717529704408
$arr[0,5] = 0.0297544130519969
$arr[0,6] = -3.55300336523739
$arr[0,7] = 0.00289091582387647
$arr[0,8] = -0.113563515145668
$arr[0,9] = 0.0336837164134808
$arr[0,10] = -3.37146631184136
$arr[0,11] = 0.00414481089215795
$arr[0,12] = -0.189329348767302
$arr[0,13] = 0.0458720525212922
$arr[0,14] = -4.1273354550556
$arr[0,15] = 0.000810087814164884
$ws.Range("B15:Q15").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = -0.0683043295252697
$arr[0,1] = 0.0108171873360851
$arr[0,2] = -6.31442605208594
$arr[0,3] = 0.0000405521819387169
$arr[0,4] = -0.0721799613511764
$arr[0,5] = 0.0136798160148092
$arr[0,6] = -5.276383927462
$arr[0,7] = 0.000201187304688376
$arr[0,8] = -0.0475684239233929
$arr[0,9] = 0.01355040163718
$arr[0,10] = -3.5104807368125
$arr[0,11] = 0.00434104559855776
$arr[0,12] = -0.089582639154158
$arr[0,13] = 0.058465668461721
$arr[0,14] = -1.53222637337688
$arr[0,15] = 0.14661159688113
$ws.Range("B16:Q16").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.000210308295272032
$arr[0,1] = 0.0000175536726668886
$arr[0,2] = 11.9808714257692
$arr[0,3] = 0.00000000439899917808996
$arr[0,4] = 0.000232237467112555
$arr[0,5] = 0.0000209572564315405
$arr[0,6] = 11.0814823434159
$arr[0,7] = 0.0000000134479151975639
$arr[0,8] = 0.000283800500778588
$arr[0,9] = 0.0000241992576746154
$arr[0,10] = 11.7276531617038
$arr[0,11] = 0.00000000522256576225549
$arr[0,12] = 0.00032244651434821
$arr[0,13] = 0.000018574674824369
$arr[0,14] = 17.3594702139914
$arr[0,15] = 0.0000000000183721069439545
$ws.Range("B17:Q17").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = -0.0142092025634642
$arr[0,1] = 0.00468718208522641
$arr[0,2] = -3.03150214886048
$arr[0,3] = 0.00823785790828971
$arr[0,4] = -0.014359789753984
$arr[0,5] = 0.00416174317579176
$arr[0,6] = -3.45042669560024
$arr[0,7] = 0.00341917723362288
$arr[0,8] = -0.0102000492509734
$arr[0,9] = 0.00494152362921783
$arr[0,10] = -2.06415065804066
$arr[0,11] = 0.0562043682100922
$arr[0,12] = -0.0159012868876992
$arr[0,13] = 0.00539908433613794
$arr[0,14] = -2.94518216380996
$arr[0,15] = 0.0100063167255513
$ws.Range("B18:Q18").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.153183395653443
$arr[0,1] = 0.0458305097436391
$arr[0,2] = 3.3423890877562
$arr[0,3] = 0.00453539464218464
$arr[0,4] = 0.159329432519212
$arr[0,5] = 0.0428326357508374
$arr[0,6] = 3.71981386917327
$arr[0,7] = 0.00205015028425931
$arr[0,8] = 0.149566735623255
$arr[0,9] = 0.0482358563262514
$arr[0,10] = 3.10073764652659
$arr[0,11] = 0.00739160496022502
$arr[0,12] = 0.199284255872885
$arr[0,13] = 0.0469375056631387
$arr[0,14] = 4.24573596439294
$arr[0,15] = 0.000714247004335106
$ws.Range("B19:Q19").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = -0.0381397192298122
$arr[0,1] = 0.015761306841851
$arr[0,2] = -2.41983229008268
$arr[0,3] = 0.0360898085786036
$arr[0,4] = -0.0491698855312385
$arr[0,5] = 0.0182217522392511
$arr[0,6] = -2.69841697360601
$arr[0,7] = 0.0223461899055611
$arr[0,8] = -0.0211932910089965
$arr[0,9] = 0.0371642160746143
$arr[0,10] = -0.570260676733957
$arr[0,11] = 0.581122957167458
$arr[0,12] = 0.00482980993526434
$arr[0,13] = 0.0403054680336961
$arr[0,14] = 0.119830141439532
$arr[0,15] = 0.907001328170302
$ws.Range("B20:Q20").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.107003625681054
$arr[0,1] = 0.00984245580516083
$arr[0,2] = 10.8716389282589
$arr[0,3] = 0.00000000785571575580565
$arr[0,4] = 0.117267572087361
$arr[0,5] = 0.0141138106426125
$arr[0,6] = 8.3087108830344
$arr[0,7] = 0.000000323415600308188
$arr[0,8] = 0.150459932752775
$arr[0,9] = 0.0172584727416505
$arr[0,10] = 8.71803287608784
$arr[0,11] = 0.000000165248355332371
$arr[0,12] = 0.150000926438644
$arr[0,13] = 0.0262051551330143
$arr[0,14] = 5.72409992145655
$arr[0,15] = 0.0000301902534546144
$ws.Range("B21:Q21").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 0.0549393783896115
$arr[0,1] = 0.00477377904546259
$arr[0,2] = 11.5085716926573
$arr[0,3] = 0.00000000218520701475582
$arr[0,4] = 0.0589559624589034
$arr[0,5] = 0.00633930546740103
$arr[0,6] = 9.30006650761287
$arr[0,7] = 0.0000000493016160813383
$arr[0,8] = 0.073359015433885
$arr[0,9] = 0.0116419226256619
$arr[0,10] = 6.30128010576036
$arr[0,11] = 0.00000842347265974575
$arr[0,12] = 0.0710564673722614
$arr[0,13] = 0.0136442949241721
$arr[0,14] = 5.20777861862093
$arr[0,15] = 0.0000736314645760792
$ws.Range("B22:Q22").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = -0.0536962940120698
$arr[0,1] = 0.00986922580440092
$arr[0,2] = -5.44078077412367
$arr[0,3] = 0.0000347453270886651
$arr[0,4] = -0.0549039593959274
$arr[0,5] = 0.00971822108264136
$arr[0,6] = -5.64958945974141
$arr[0,7] = 0.0000223374778883745
$arr[0,8] = -0.0544685616914582
$arr[0,9] = 0.015823049050442
$arr[0,10] = -3.44235561160298
$arr[0,11] = 0.002869292364151
$arr[0,12] = -0.0280449898694483
$arr[0,13] = 0.0157819743313352
$arr[0,14] = -1.77702670658669
$arr[0,15] = 0.0922721925553008
$ws.Range("B23:Q23").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.0267473950661548
$arr[0,1] = 0.0121290869942179
$arr[0,2] = 2.20522740738076
$arr[0,3] = 0.0440136945101541
$arr[0,4] = 0.0274013673122529
$arr[0,5] = 0.0136130515115439
$arr[0,6] = 2.01287472459915
$arr[0,7] = 0.0630682730820312
$arr[0,8] = 0.050045035526746
$arr[0,9] = 0.0207247972834977
$arr[0,10] = 2.4147418593375
$arr[0,11] = 0.0294524403119403
$ws.Range("B24:M24").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.0727690058715467
$arr[0,1] = 0.023544461512789
$arr[0,2] = -3.09070588987646
$arr[0,3] = 0.0202653094440712
$arr[0,4] = -0.068944093898408
$arr[0,5] = 0.0293571531985947
$arr[0,6] = -2.34845979213367
$arr[0,7] = 0.0554795912753991
$arr[0,8] = -0.0647711336198183
$arr[0,9] = 0.0356167489704901
$arr[0,10] = -1.81855827643011
$arr[0,11] = 0.116871094765824
$ws.Range("B25:M25").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.179417719312666
$arr[0,1] = 0.0413535916762189
$arr[0,2] = -4.3386248216946
$arr[0,3] = 0.000815368565572474
$arr[0,4] = -0.184610696156454
$arr[0,5] = 0.0414953584830741
$arr[0,6] = -4.44894809697227
$arr[0,7] = 0.000696031197702425
$arr[0,8] = -0.143463591030625
$arr[0,9] = 0.048689810429095
$arr[0,10] = -2.9464807886148
$arr[0,11] = 0.0115772539691394
$ws.Range("B26:M26").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.136801817116175
$arr[0,1] = 0.0782803690864455
$arr[0,2] = -1.74758778877376
$arr[0,3] = 0.110689545488128
$arr[0,4] = -0.117710952355572
$arr[0,5] = 0.0819689646613442
$arr[0,6] = -1.43604293212555
$arr[0,7] = 0.18098866052077
$arr[0,8] = -0.0755130863571989
$arr[0,9] = 0.0855007428531531
$arr[0,10] = -0.883186319057977
$arr[0,11] = 0.397088850439362
$ws.Range("B27:M27").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.133978745474261
$arr[0,1] = 0.0354521476887659
$arr[0,2] = -3.77914327364478
$arr[0,3] = 0.00482523833248853
$arr[0,4] = -0.124722083179914
$arr[0,5] = 0.0331067739016527
$arr[0,6] = -3.76726779692925
$arr[0,7] = 0.00504205218937793
$arr[0,8] = -0.129637044266528
$arr[0,9] = 0.0368336451003934
$arr[0,10] = -3.51952797267798
$arr[0,11] = 0.00730456331181071
$ws.Range("B28:M28").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.141706757562063
$arr[0,1] = 0.201290737260586
$arr[0,2] = -0.703990454258276
$arr[0,3] = 0.571921102362629
$arr[0,4] = -0.110856257039908
$arr[0,5] = 0.207644072051943
$arr[0,6] = -0.533876339181871
$arr[0,7] = 0.659917955659709
$arr[0,8] = -0.0976900284260015
$arr[0,9] = 0.210291818750349
$arr[0,10] = -0.464545073633965
$arr[0,11] = 0.699302253729195
$ws.Range("B29:M29").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.100207680598556
$arr[0,1] = 0.0241298060953741
$arr[0,2] = -4.15285892486995
$arr[0,3] = 0.000833693648638827
$arr[0,4] = -0.124176382196076
$arr[0,5] = 0.0302808088704361
$arr[0,6] = -4.10082777931708
$arr[0,7] = 0.000923419373718214
$arr[0,8] = -0.12898233903257
$arr[0,9] = 0.0383737451204796
$arr[0,10] = -3.36121321042844
$arr[0,11] = 0.00423699057070018
$ws.Range("B30:M30").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.074397128684255
$arr[0,1] = 0.0260842576718059
$arr[0,2] = -2.85218500830368
$arr[0,3] = 0.0137750432730592
$arr[0,4] = -0.0834994657792398
$arr[0,5] = 0.0254574283812728
$arr[0,6] = -3.2799646739128
$arr[0,7] = 0.00608905118832989
$arr[0,8] = -0.11724525249082
$arr[0,9] = 0.0359509046366676
$arr[0,10] = -3.26126014562752
$arr[0,11] = 0.00630123776895843
$ws.Range("B31:M31").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.0179086370568567
$arr[0,1] = 0.12078119581384
$arr[0,2] = -0.148273387559925
$arr[0,3] = 0.894242171415072
$arr[0,4] = -0.0386041040131227
$arr[0,5] = 0.1206329561241
$arr[0,6] = -0.320012915653074
$arr[0,7] = 0.776240700788529
$arr[0,8] = -0.0249629268394769
$arr[0,9] = 0.0500788984421327
$arr[0,10] = -0.498471963562101
$arr[0,11] = 0.669692083141968
$ws.Range("B32:M32").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.0447473477007656
$arr[0,1] = 0.0563403721249598
$arr[0,2] = -0.794232377477353
$arr[0,3] = 0.450551364624079
$arr[0,4] = -0.0891147649458146
$arr[0,5] = 0.057546185870721
$arr[0,6] = -1.54857813072118
$arr[0,7] = 0.161997077668497
$arr[0,8] = -0.0392849418827011
$arr[0,9] = 0.0652744994774341
$arr[0,10] = -0.601842100624337
$arr[0,11] = 0.565213598566131
$ws.Range("B33:M33").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.094270009954612
$arr[0,1] = 0.0233485928941462
$arr[0,2] = -4.03750283291147
$arr[0,3] = 0.00574173756346435
$arr[0,4] = -0.0658938776653298
$arr[0,5] = 0.0307046502231931
$arr[0,6] = -2.14605531039582
$arr[0,7] = 0.071910332274278
$arr[0,8] = -0.0889420587555716
$arr[0,9] = 0.0477034389976704
$arr[0,10] = -1.86447896890442
$arr[0,11] = 0.108239193960291
$ws.Range("B34:M34").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.116051887583236
$arr[0,1] = 0.0463706807043698
$arr[0,2] = -2.50269967618352
$arr[0,3] = 0.0239254901169165
$arr[0,4] = -0.141785441268686
$arr[0,5] = 0.0442249366953392
$arr[0,6] = -3.20600665288525
$arr[0,7] = 0.00551743349082206
$arr[0,8] = -0.116753997419382
$arr[0,9] = 0.0478923699487857
$arr[0,10] = -2.43784130007001
$arr[0,11] = 0.0277497769217896
$ws.Range("B35:M35").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.00188345262812793
$arr[0,1] = 0.0645654673330066
$arr[0,2] = 0.0291712072401448
$arr[0,3] = 0.977529074107264
$arr[0,4] = -0.0355550781861411
$arr[0,5] = 0.0881886579576326
$arr[0,6] = -0.403170645858137
$arr[0,7] = 0.698813290185855
$arr[0,8] = 0.0554630192755981
$arr[0,9] = 0.0894297800316212
$arr[0,10] = 0.620185124641782
$arr[0,11] = 0.555828049980281
$ws.Range("B36:M36").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = -0.0391248247693395
$arr[0,1] = 0.150605150990689
$arr[0,2] = -0.259784107727885
$arr[0,3] = 0.799545007554951
$arr[0,4] = 0.0328207381819786
$arr[0,5] = 0.144429160719982
$arr[0,6] = 0.227244539941703
$arr[0,7] = 0.82441930413716
$arr[0,8] = 0.125738107830748
$arr[0,9] = 0.152602443720071
$arr[0,10] = 0.823958678285633
$arr[0,11] = 0.428348956407248
$ws.Range("B37:M37").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.242422791760949
$arr[0,1] = 0.149465690418441
$arr[0,2] = 1.62192936106118
$arr[0,3] = 0.210535815949125
$arr[0,4] = 0.523770124510843
$arr[0,5] = 0.200716071282856
$arr[0,6] = 2.60950765508322
$arr[0,7] = 0.0964224962151562
$arr[0,8] = 0.603574755170635
$arr[0,9] = 0.140378063562112
$arr[0,10] = 4.29963727846677
$arr[0,11] = 0.0350454329379503
$ws.Range("B38:M38").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.408607328184208
$arr[0,1] = 0.0347070762370814
$arr[0,2] = 11.7730264973357
$arr[0,3] = 0.00362535579887727
$arr[0,4] = 0.35338523065181
$arr[0,5] = 0.055171579119126
$arr[0,6] = 6.4052042064771
$arr[0,7] = 0.016969511659237
$arr[0,8] = 0.357038416094967
$arr[0,9] = 0.0636276671842467
$arr[0,10] = 5.61137052315765
$arr[0,11] = 0.0225622068693846
$ws.Range("B39:M39").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.630122782672917
$arr[0,1] = 0.0518648656744356
$arr[0,2] = 12.1493187050421
$arr[0,3] = 0.00051060462823688
$arr[0,4] = 0.587193595365284
$arr[0,5] = 0.0754637958844846
$arr[0,6] = 7.78112985813918
$arr[0,7] = 0.002927273032098
$arr[0,8] = 0.607972988584232
$arr[0,9] = 0.0879885304057641
$arr[0,10] = 6.90968454388919
$arr[0,11] = 0.00442621939114473
$ws.Range("B40:M40").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.638549155045028
$arr[0,1] = 0.0621407404934607
$arr[0,2] = 10.2758536505085
$arr[0,3] = 0.00003249265486439
$arr[0,4] = 0.588759306198064
$arr[0,5] = 0.0890883794296028
$arr[0,6] = 6.60871047343833
$arr[0,7] = 0.000507284513743361
$arr[0,8] = 0.600370123632559
$arr[0,9] = 0.101082038505468
$arr[0,10] = 5.93943427051172
$arr[0,11] = 0.000920241977739678
$ws.Range("B41:M41").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.657537729911562
$arr[0,1] = 0.0491342758961082
$arr[0,2] = 13.3824650494879
$arr[0,3] = 0.0000235163172896176
$arr[0,4] = 0.615825940605714
$arr[0,5] = 0.0761238151358747
$arr[0,6] = 8.08979344383247
$arr[0,7] = 0.00039343887426625
$arr[0,8] = 0.620867468242027
$arr[0,9] = 0.0895222220843181
$arr[0,10] = 6.93534469751267
$arr[0,11] = 0.000860224330447681
$ws.Range("B42:M42").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.605309079423922
$arr[0,1] = 0.0482941192046902
$arr[0,2] = 12.5338051380205
$arr[0,3] = 0.000000741119755418709
$arr[0,4] = 0.547982584458141
$arr[0,5] = 0.0829481546888635
$arr[0,6] = 6.60632640368685
$arr[0,7] = 0.000128486697956472
$arr[0,8] = 0.546850255560341
$arr[0,9] = 0.0983302461454261
$arr[0,10] = 5.56136363933814
$arr[0,11] = 0.000443481727493624
$ws.Range("B43:M43").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.00000403887212738514
$arr[0,1] = 0.000142473491862056
$arr[0,2] = 0.0283482356935255
$arr[0,3] = 0.977999713040412
$arr[0,4] = 0.000142636162499181
$arr[0,5] = 0.000150858511177591
$arr[0,6] = 0.945496289110728
$arr[0,7] = 0.369557091617099
$arr[0,8] = -0.0001245690876513
$arr[0,9] = 0.000253389372062752
$arr[0,10] = -0.491611335697417
$arr[0,11] = 0.635317453582563
$ws.Range("B44:M44").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.00271290401391287
$arr[0,1] = 0.00442881339842759
$arr[0,2] = 0.612557759799964
$arr[0,3] = 0.548527513840926
$arr[0,4] = -0.00242431006387328
$arr[0,5] = 0.00445354219172934
$arr[0,6] = -0.544355472454141
$arr[0,7] = 0.593490633266188
$arr[0,8] = -0.0109708542650427
$arr[0,9] = 0.00613784936763943
$arr[0,10] = -1.78741015100246
$arr[0,11] = 0.092255128786812
$ws.Range("B45:M45").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.00334686560434135
$arr[0,1] = 0.00233019613870859
$arr[0,2] = 1.43630209867063
$arr[0,3] = 0.238001057403457
$arr[0,4] = 0.0105181007010106
$arr[0,5] = 0.00385904361565796
$arr[0,6] = 2.72557186405817
$arr[0,7] = 0.0654357758670905
$arr[0,8] = 0.0153130593443016
$arr[0,9] = 0.00231693757292333
$arr[0,10] = 6.60918080972755
$arr[0,11] = 0.00523998834815388
$ws.Range("B46:M46").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.00517500591611341
$arr[0,1] = 0.00710709184081542
$arr[0,2] = -0.728146762701699
$arr[0,3] = 0.484602616030662
$arr[0,4] = -0.00874858193967026
$arr[0,5] = 0.0095633436943417
$arr[0,6] = -0.914803673201298
$arr[0,7] = 0.383078982261738
$ws.Range("B47:I47").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.055647141281049
$arr[0,1] = 0.0322362382512921
$arr[0,2] = -1.72622937103459
$arr[0,3] = 0.107174599422883
$arr[0,4] = -0.0355836369677184
$arr[0,5] = 0.0372703368607687
$arr[0,6] = -0.954744173648032
$arr[0,7] = 0.356603886560485
$ws.Range("B48:I48").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.00850083333593829
$arr[0,1] = 0.00835299340666677
$arr[0,2] = -1.01769903579159
$arr[0,3] = 0.332303711806271
$arr[0,4] = -0.00585680165835267
$arr[0,5] = 0.00954442772168712
$arr[0,6] = -0.613635707570469
$arr[0,7] = 0.554379207308233
$ws.Range("B49:I49").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.00702466463209921
$arr[0,1] = 0.00626785703607644
$arr[0,2] = -1.12074423390112
$arr[0,3] = 0.280897027094083
$arr[0,4] = -0.00535520119284443
$arr[0,5] = 0.00817919942165448
$arr[0,6] = -0.654734151445985
$arr[0,7] = 0.523308880340958
$ws.Range("B50:I50").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.00127323073443998
$arr[0,1] = 0.0015229270612078
$arr[0,2] = -0.836041834748285
$arr[0,3] = 0.547339568165121
$arr[0,4] = -0.000655052261136138
$arr[0,5] = 0.00560251367093152
$arr[0,6] = -0.116921135692155
$arr[0,7] = 0.924698638286826
$ws.Range("B51:I51").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.0184050510467514
$arr[0,1] = 0.0132837296306737
$arr[0,2] = 1.38553339750697
$arr[0,3] = 0.190452604570735
$arr[0,4] = 0.0230650381452688
$arr[0,5] = 0.0145342555719752
$arr[0,6] = 1.58694320676063
$arr[0,7] = 0.137215085200518
$ws.Range("B52:I52").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.00875103135470041
$arr[0,1] = 0.0121348514421548
$arr[0,2] = 0.721148618622604
$arr[0,3] = 0.4819143972567
$arr[0,4] = 0.0155424713472888
$arr[0,5] = 0.0171507559172309
$arr[0,6] = 0.906226607287537
$arr[0,7] = 0.379046718689801
$ws.Range("B53:I53").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.00000476967739594562
$arr[0,1] = 0.00000547899822054111
$arr[0,2] = 0.870538226872168
$arr[0,3] = 0.40127638689801
$arr[0,4] = 0.00000453032124940568
$arr[0,5] = 0.00000650627599133779
$arr[0,6] = 0.696300196216266
$arr[0,7] = 0.500821721192307
$ws.Range("B54:I54").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.0000129676315070935
$arr[0,1] = 0.00000271680297708411
$arr[0,2] = 4.77312179663887
$arr[0,3] = 0.000238361870474816
$arr[0,4] = 0.0000145470271275458
$arr[0,5] = 0.00000359086014389871
$arr[0,6] = 4.05112606578758
$arr[0,7] = 0.00120428670711552
$ws.Range("B55:I55").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = -0.481449696900955
$arr[0,1] = 0.173178235912867
$arr[0,2] = -2.780082002586
$arr[0,3] = 0.0277532294095023
$ws.Range("B56:E56").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = -0.140355165345535
$arr[0,1] = 0.0601237665100454
$arr[0,2] = -2.33443733639151
$arr[0,3] = 0.0366297114791482
$ws.Range("B57:E57").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = -0.396019620539022
$arr[0,1] = 0.146374970730123
$arr[0,2] = -2.70551460105278
$arr[0,3] = 0.0165350548521122
$ws.Range("B58:E58").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 0.302422210778769
$arr[0,1] = 0.0858985685128793
$arr[0,2] = 3.52068976252409
$arr[0,3] = 0.00301811913782489
$ws.Range("B59:E59").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = -0.665033319291968
$arr[0,1] = 0.0926845388943302
$arr[0,2] = -7.17523469637341
$arr[0,3] = 0.000149617202537518
$ws.Range("B60:E60").Value = $arr
